$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 27965
$ws.Range("E2").Value = -292
$ws.Range("F2").Value = -292
$ws.Range("G2").Value = 1154
$ws.Range("H2").Value = 1006
$ws.Range("I2").Value = 1281
$ws.Range("J2").Value = -275
$ws.Range("K2").Value = 38954
$ws.Range("L2").Value = 9831
$ws.Range("M2").Value = 29122
$ws.Range("N2").Value = 24487
$ws.Range("O2").Value = 4635
$ws.Range("P2").Value = 92
$ws.Range("Q2").Value = 2127
$ws.Range("R2").Value = -1664
$ws.Range("S2").Value = -113
$ws.Range("T2").Value = 1743
$ws.Range("U2").Value = 384
$ws.Range("V2").Value = 2947
$ws.Range("W2").Value = -1.05
$ws.Range("X2").Value = 3.6
$ws.Range("Y2").Value = 5.37
$ws.Range("Z2").Value = 2.57
$ws.Range("AA2").Value = 33.76
$ws.Range("AB2").Value = 26708.9
$ws.Range("AC2").Value = 69540
$ws.Range("AD2").Value = 17.8
$ws.Range("AE2").Value = 1423568
$ws.Range("AF2").Value = 0.87
$ws.Range("AG2").Value = 7500
$ws.Range("AH2").Value = 0.61
$ws.Range("AI2").Value = 10.07
$ws.Range("AJ2").Value = 1842040

$ws.Range("D3").Value = 26154
$ws.Range("E3").Value = -55
$ws.Range("F3").Value = -55
$ws.Range("G3").Value = 1147
$ws.Range("H3").Value = 899
$ws.Range("I3").Value = 1289
$ws.Range("J3").Value = -390
$ws.Range("K3").Value = 39879
$ws.Range("L3").Value = 9844
$ws.Range("M3").Value = 30034
$ws.Range("N3").Value = 25805
$ws.Range("O3").Value = 4230
$ws.Range("P3").Value = 92
$ws.Range("Q3").Value = 2101
$ws.Range("R3").Value = -1291
$ws.Range("S3").Value = -340
$ws.Range("T3").Value = 1041
$ws.Range("U3").Value = 1060
$ws.Range("V3").Value = 2614
$ws.Range("W3").Value = -0.21
$ws.Range("X3").Value = 3.44
$ws.Range("Y3").Value = 5.13
$ws.Range("Z3").Value = 2.28
$ws.Range("AA3").Value = 32.78
$ws.Range("AB3").Value = 27949.65
$ws.Range("AC3").Value = 69989
$ws.Range("AD3").Value = 15.67
$ws.Range("AE3").Value = 1500147
$ws.Range("AF3").Value = 0.73
$ws.Range("AG3").Value = 9750
$ws.Range("AH3").Value = 0.89
$ws.Range("AI3").Value = 13.01
$ws.Range("AJ3").Value = 1842040

$ws.Range("D4").Value = 26541
$ws.Range("E4").Value = -42
$ws.Range("F4").Value = -42
$ws.Range("G4").Value = 1715
$ws.Range("H4").Value = 1274
$ws.Range("I4").Value = 1624
$ws.Range("J4").Value = -350
$ws.Range("K4").Value = 41349
$ws.Range("L4").Value = 9926
$ws.Range("M4").Value = 31423
$ws.Range("N4").Value = 27231
$ws.Range("O4").Value = 4192
$ws.Range("P4").Value = 92
$ws.Range("Q4").Value = 794
$ws.Range("R4").Value = 126
$ws.Range("S4").Value = -534
$ws.Range("T4").Value = 809
$ws.Range("U4").Value = -15
$ws.Range("V4").Value = 1829
$ws.Range("W4").Value = -0.16
$ws.Range("X4").Value = 4.8
$ws.Range("Y4").Value = 6.12
$ws.Range("Z4").Value = 3.14
$ws.Range("AA4").Value = 31.59
$ws.Range("AB4").Value = 29503.1
$ws.Range("AC4").Value = 88137
$ws.Range("AD4").Value = 12.07
$ws.Range("AE4").Value = 1583095
$ws.Range("AF4").Value = 0.67
$ws.Range("AG4").Value = 10000
$ws.Range("AH4").Value = 0.9399999999999999
$ws.Range("AI4").Value = 10.6
$ws.Range("AJ4").Value = 1842040

$ws.Range("D5").Value = 37249
$ws.Range("E5").Value = 1594
$ws.Range("F5").Value = 1594
$ws.Range("G5").Value = 3346
$ws.Range("H5").Value = 2723
$ws.Range("I5").Value = 2446
$ws.Range("J5").Value = 277
$ws.Range("K5").Value = 48037
$ws.Range("L5").Value = 13342
$ws.Range("M5").Value = 34695
$ws.Range("N5").Value = 29658
$ws.Range("O5").Value = 5037
$ws.Range("P5").Value = 92
$ws.Range("Q5").Value = 1838
$ws.Range("R5").Value = -3786
$ws.Range("S5").Value = 1809
$ws.Range("T5").Value = 4287
$ws.Range("U5").Value = -2449
$ws.Range("V5").Value = 3139
$ws.Range("W5").Value = 4.28
$ws.Range("X5").Value = 7.31
$ws.Range("Y5").Value = 8.6
$ws.Range("Z5").Value = 6.09
$ws.Range("AA5").Value = 38.46
$ws.Range("AB5").Value = 31994.86
$ws.Range("AC5").Value = 132786
$ws.Range("AD5").Value = 8.09
$ws.Range("AE5").Value = 1724178
$ws.Range("AF5").Value = 0.62
$ws.Range("AG5").Value = 10000
$ws.Range("AH5").Value = 0.93
$ws.Range("AI5").Value = 7.03
$ws.Range("AJ5").Value = 1842040

$ws.Range("D6").Value = 29714
$ws.Range("E6").Value = -1089
$ws.Range("F6").Value = -1089
$ws.Range("G6").Value = 667
$ws.Range("H6").Value = 462
$ws.Range("I6").Value = 1107
$ws.Range("K6").Value = 47216
$ws.Range("L6").Value = 12142
$ws.Range("M6").Value = 35074
$ws.Range("N6").Value = 30412
$ws.Range("P6").Value = 92
$ws.Range("Q6").Value = 1190
$ws.Range("R6").Value = -1757
$ws.Range("S6").Value = 550
$ws.Range("T6").Value = 1126
$ws.Range("U6").Value = 65
$ws.Range("V6").Value = 2943
$ws.Range("W6").Value = -3.67
$ws.Range("X6").Value = 1.56
$ws.Range("Y6").Value = 3.69
$ws.Range("Z6").Value = 0.97
$ws.Range("AA6").Value = 34.62
$ws.Range("AB6").Value = 33132.32
$ws.Range("AC6").Value = 60108
$ws.Range("AD6").Value = 12.39
$ws.Range("AE6").Value = 1768014
$ws.Range("AF6").Value = 0.42
$ws.Range("AG6").Value = 10000
$ws.Range("AH6").Value = 1.34
$ws.Range("AI6").Value = 15.54
$ws.Range("AJ6").Value = 1842040

$ws.Range("D7").Value = 32140
$ws.Range("E7").Value = 944
$ws.Range("G7").Value = 2904
$ws.Range("H7").Value = 2469
$ws.Range("I7").Value = 2493
$ws.Range("K7").Value = 49620
$ws.Range("L7").Value = 12400
$ws.Range("M7").Value = 37220
$ws.Range("N7").Value = 32580
$ws.Range("P7").Value = 90
$ws.Range("Q7").Value = 2160
$ws.Range("R7").Value = -810
$ws.Range("S7").Value = 110
$ws.Range("T7").Value = 800
$ws.Range("U7").Value = 1120
$ws.Range("W7").Value = 2.94
$ws.Range("X7").Value = 7.68
$ws.Range("Y7").Value = 7.92
$ws.Range("Z7").Value = 5.1
$ws.Range("AA7").Value = 33.32
$ws.Range("AC7").Value = 135339
$ws.Range("AD7").Value = 4.63
$ws.Range("AE7").Value = 1894038
$ws.Range("AF7").Value = 0.33
$ws.Range("AG7").Value = 10000
$ws.Range("AH7").Value = 1.6
$ws.Range("AI7").Value = 7.39

$ws.Range("D8").Value = 33210
$ws.Range("E8").Value = 1192
$ws.Range("G8").Value = 3351
$ws.Range("H8").Value = 2681
$ws.Range("I8").Value = 2600
$ws.Range("K8").Value = 52080
$ws.Range("L8").Value = 12510
$ws.Range("M8").Value = 39570
$ws.Range("N8").Value = 34850
$ws.Range("P8").Value = 90
$ws.Range("Q8").Value = 2290
$ws.Range("R8").Value = -1010
$ws.Range("S8").Value = 110
$ws.Range("T8").Value = 1000
$ws.Range("U8").Value = 1170
$ws.Range("W8").Value = 3.59
$ws.Range("X8").Value = 8.07
$ws.Range("Y8").Value = 7.71
$ws.Range("Z8").Value = 5.27
$ws.Range("AA8").Value = 31.61
$ws.Range("AC8").Value = 141148
$ws.Range("AD8").Value = 4.44
$ws.Range("AE8").Value = 2026005
$ws.Range("AF8").Value = 0.31
$ws.Range("AG8").Value = 10000
$ws.Range("AH8").Value = 1.6
$ws.Range("AI8").Value = 7.08

$ws.Range("D9").Value = 33460
$ws.Range("E9").Value = 1220
$ws.Range("G9").Value = 3390
$ws.Range("H9").Value = 2710
$ws.Range("I9").Value = 2630
$ws.Range("K9").Value = 54500
$ws.Range("L9").Value = 12540
$ws.Range("M9").Value = 41960
$ws.Range("N9").Value = 37160
$ws.Range("P9").Value = 90
$ws.Range("Q9").Value = 2390
$ws.Range("R9").Value = -1010
$ws.Range("S9").Value = 110
$ws.Range("T9").Value = 1000
$ws.Range("U9").Value = 1260
$ws.Range("W9").Value = 3.65
$ws.Range("X9").Value = 8.1
$ws.Range("Y9").Value = 7.31
$ws.Range("Z9").Value = 5.08
$ws.Range("AA9").Value = 29.89
$ws.Range("AC9").Value = 142776
$ws.Range("AD9").Value = 4.38
$ws.Range("AE9").Value = 2160297
$ws.Range("AF9").Value = 0.29
$ws.Range("AG9").Value = 10000
$ws.Range("AH9").Value = 1.6
$ws.Range("AI9").Value = 7
